# "Added plot aggregation to year and parameter set exporting"
#
# Observed edits in this commit are all workbook/sheet *view* state plus a
# worksheet rename:
#   - "crop sheet lookup" worksheet renamed to "crop_sheet_lookup"
#   - the active/selected tab moves from "simulation_control" (first sheet)
#     to "crop_sheet_lookup" (last sheet), whose selected cell becomes L16
#   - the host OS window's saved position/size (xWindow/yWindow/
#     windowWidth/windowHeight) changed too, but that is pure window-chrome
#     state captured by Excel from the desktop window manager, not
#     something exposed on the Application/Workbook/Window object model -
#     there is nothing in this headless session to move/resize, so it is
#     left as-is.

$wb = $excel.ActiveWorkbook

# Rename "crop sheet lookup" -> "crop_sheet_lookup"
$lookupSheet = $wb.Worksheets.Item("crop sheet lookup")
$lookupSheet.Name = "crop_sheet_lookup"

# Activate it (this both sets bookViews/workbookView activeTab and moves
# tabSelected="1" off of whichever sheet previously had it) and place the
# selection at L16.
$lookupSheet.Select()
$lookupSheet.Range("L16").Select()
